# Generate Report for Archive
# The file "7954c250-02cb-4228-8d74-e597c501fcaa.md" moved from "Ready for
# handoff" to "In Translation" status. Update the Status column on every
# sheet that tracks it.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B5").Value = "In Translation"
$overview.Range("C5").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B5").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B5").Value = "In Translation"
